$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.888.39"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "2.284.52"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.15"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("E7").Value = "  -0.67%  "
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.59"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.991"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.19"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "2.631.06"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "2.283.55"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "42.779.55"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.41%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +24.30%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000105"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.73"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.52"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.01"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.36%  "
$ws.Range("E25").Value = "  -3.11%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.14"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +21.72%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.45"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.93"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.39"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0870"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.59"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("E36").Value = "  -3.89%  "
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.77"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("E40").Value = "  -3.90%  "
$ws.Range("E41").Value = "  +5.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.12"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.230"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.77"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.13"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.39%  "
$ws.Range("D47").Value = "1.731.03"
$ws.Range("E47").Value = "  +8.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "113.34"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "78.52"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.79%  "
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.18"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.16%  "
